$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E; existing D:K data shifts to F:M.
$ws.Columns("D:E").Insert()

# The freshly inserted D:E columns don't inherit the number formats that
# used to live in (now-shifted) column F onward, so pull them back in -
# this matches Excel's "insert copied cells" behavior for this edit.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# New quarter data (D) and prior quarter (E) for every populated row.
$newQuarterData = @(
    @{Row=7; D=43465; E=43373},
    @{Row=8; D=11400; E=11000},
    @{Row=9; D="NA"; E="NA"},
    @{Row=10; D="NA"; E="NA"},
    @{Row=11; D=""; E=""},
    @{Row=12; D="NA"; E="NA"},
    @{Row=13; D=0; E=0},
    @{Row=14; D=0; E=0},
    @{Row=15; D=0; E=0},
    @{Row=16; D=""; E=""},
    @{Row=17; D=1800; E=1700},
    @{Row=18; D=9600; E=9300},
    @{Row=19; D=""; E=""},
    @{Row=20; D=-4100; E=-3900},
    @{Row=21; D=5700; E=5500},
    @{Row=22; D=0; E=0},
    @{Row=23; D=5500; E=5400},
    @{Row=24; D=1000; E=1000},
    @{Row=25; D=0; E=0},
    @{Row=26; D=4600; E=4400},
    @{Row=27; D=4600; E=4400},
    @{Row=28; D=0; E=0},
    @{Row=29; D="NA"; E="NA"},
    @{Row=30; D=0; E=0},
    @{Row=31; D=0; E=0},
    @{Row=32; D=4100; E=3900},
    @{Row=33; D=4600; E=4400},
    @{Row=34; D=0; E=0},
    @{Row=35; D=4600; E=4400},
    @{Row=38; D=43465; E=43373},
    @{Row=39; D=""; E=""},
    @{Row=40; D=""; E=""},
    @{Row=41; D=31900; E=35100},
    @{Row=42; D=6500; E=6800},
    @{Row=43; D=0; E=0},
    @{Row=44; D=0; E=0},
    @{Row=45; D=0; E=0},
    @{Row=46; D=0; E=0},
    @{Row=47; D=0; E=0},
    @{Row=48; D=13300; E=13300},
    @{Row=49; D=0; E=0},
    @{Row=50; D=0; E=0},
    @{Row=51; D=0; E=0},
    @{Row=52; D=0; E=0},
    @{Row=53; D=0; E=0},
    @{Row=54; D=1195300; E=1187900},
    @{Row=55; D=""; E=""},
    @{Row=56; D=""; E=""},
    @{Row=57; D=0; E=0},
    @{Row=58; D=0; E=0},
    @{Row=59; D=11900; E=9800},
    @{Row=60; D=0; E=0},
    @{Row=61; D=0; E=0},
    @{Row=62; D=0; E=0},
    @{Row=63; D=0; E=0},
    @{Row=64; D=0; E=0},
    @{Row=65; D=0; E=0},
    @{Row=66; D=1091200; E=1088200},
    @{Row=67; D=""; E=""},
    @{Row=68; D=0; E=0},
    @{Row=69; D=0; E=0},
    @{Row=70; D=0; E=0},
    @{Row=71; D=0; E=0},
    @{Row=72; D=94000; E=89900},
    @{Row=73; D=0; E=0},
    @{Row=74; D=0; E=0},
    @{Row=75; D=0; E=0},
    @{Row=76; D=104100; E=99600},
    @{Row=77; D=0; E=0},
    @{Row=80; D=43465; E=43373},
    @{Row=81; D=4600; E=4400},
    @{Row=82; D=""; E=""},
    @{Row=83; D=200; E=200},
    @{Row=84; D=0; E=0},
    @{Row=85; D=0; E=0},
    @{Row=86; D=0; E=0},
    @{Row=87; D=0; E=0},
    @{Row=88; D=0; E=0},
    @{Row=89; D=6000; E=2900},
    @{Row=90; D=""; E=""},
    @{Row=91; D=-200; E=-100},
    @{Row=92; D=0; E=0},
    @{Row=93; D=0; E=0},
    @{Row=94; D=-9700; E=-28400},
    @{Row=95; D=""; E=""},
    @{Row=96; D=-400; E=-900},
    @{Row=97; D=0; E=0},
    @{Row=98; D=0; E=0},
    @{Row=99; D=0; E=0},
    @{Row=100; D=500; E=34200},
    @{Row=101; D=0; E=0},
    @{Row=102; D=-3200; E=8600}
)

foreach ($entry in $newQuarterData) {
    $ws.Cells.Item($entry.Row, 4).Value = $entry.D   # column D
    $ws.Cells.Item($entry.Row, 5).Value = $entry.E   # column E
}
